$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.00841100000001
$ws.Range("H2").Value = 75.02523300000001
$ws.Range("I2").Value = 0.4156829172908309
$ws.Range("J2").Value = 0.415682917290831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 29.546424
$ws.Range("N2").Value = 88.63927200000001
$ws.Range("O2").Value = 0.9033225104610835
$ws.Range("P2").Value = 0.9033225104610834
$ws.Range("Q2").Value = 738.9091149722642
$ws.Range("R2").Value = 6650.182034750378
$ws.Range("S2").Value = 0.3754957364029403
$ws.Range("T2").Value = 0.3754957364029403

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.00841100000001
$ws.Range("H3").Value = 75.02523300000001
$ws.Range("I3").Value = 0.4156829172908309
$ws.Range("J3").Value = 0.415682917290831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.672785333333334
$ws.Range("N3").Value = 5.018356000000001
$ws.Range("O3").Value = 0.05114204841740398
$ws.Range("P3").Value = 0.05114204841740398
$ws.Range("Q3").Value = 41.83370313077202
$ws.Range("R3").Value = 376.5033281769481
$ws.Range("S3").Value = 0.02125887588237541
$ws.Range("T3").Value = 0.02125887588237541

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.00841100000001
$ws.Range("H4").Value = 75.02523300000001
$ws.Range("I4").Value = 0.4156829172908309
$ws.Range("J4").Value = 0.415682917290831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.489401
$ws.Range("N4").Value = 4.468203
$ws.Range("O4").Value = 0.04553544112151264
$ws.Range("P4").Value = 0.04553544112151264
$ws.Range("Q4").Value = 37.24755235181101
$ws.Range("R4").Value = 335.2279711662991
$ws.Range("S4").Value = 0.01892830500551524
$ws.Range("T4").Value = 0.01892830500551524

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.496019
$ws.Range("H5").Value = 61.488057
$ws.Range("I5").Value = 0.340679180727168
$ws.Range("J5").Value = 0.3406791807271681
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.546424
$ws.Range("N5").Value = 88.63927200000001
$ws.Range("O5").Value = 0.9033225104610835
$ws.Range("P5").Value = 0.9033225104610834
$ws.Range("Q5").Value = 605.5840676860561
$ws.Range("R5").Value = 5450.256609174504
$ws.Range("S5").Value = 0.3077431727962906
$ws.Range("T5").Value = 0.3077431727962906

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.496019
$ws.Range("H6").Value = 61.488057
$ws.Range("I6").Value = 0.340679180727168
$ws.Range("J6").Value = 0.3406791807271681
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.672785333333334
$ws.Range("N6").Value = 5.018356000000001
$ws.Range("O6").Value = 0.05114204841740398
$ws.Range("P6").Value = 0.05114204841740398
$ws.Range("Q6").Value = 34.28543997492134
$ws.Range("R6").Value = 308.568959774292
$ws.Range("S6").Value = 0.01742303115555035
$ws.Range("T6").Value = 0.01742303115555035

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.496019
$ws.Range("H7").Value = 61.488057
$ws.Range("I7").Value = 0.340679180727168
$ws.Range("J7").Value = 0.3406791807271681
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.489401
$ws.Range("N7").Value = 4.468203
$ws.Range("O7").Value = 0.04553544112151264
$ws.Range("P7").Value = 0.04553544112151264
$ws.Range("Q7").Value = 30.526791194619
$ws.Range("R7").Value = 274.741120751571
$ws.Range("S7").Value = 0.01551297677532712
$ws.Range("T7").Value = 0.01551297677532712

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.65779933333333
$ws.Range("H8").Value = 43.973398
$ws.Range("I8").Value = 0.243637901982001
$ws.Range("J8").Value = 0.243637901982001
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.546424
$ws.Range("N8").Value = 88.63927200000001
$ws.Range("O8").Value = 0.9033225104610835
$ws.Range("P8").Value = 0.9033225104610834
$ws.Range("Q8").Value = 433.0855540095841
$ws.Range("R8").Value = 3897.769986086257
$ws.Range("S8").Value = 0.2200836012618525
$ws.Range("T8").Value = 0.2200836012618525

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.65779933333333
$ws.Range("H9").Value = 43.973398
$ws.Range("I9").Value = 0.243637901982001
$ws.Range("J9").Value = 0.243637901982001
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.672785333333334
$ws.Range("N9").Value = 5.018356000000001
$ws.Range("O9").Value = 0.05114204841740398
$ws.Range("P9").Value = 0.05114204841740398
$ws.Range("Q9").Value = 24.51935174374312
$ws.Range("R9").Value = 220.6741656936881
$ws.Range("S9").Value = 0.01246014137947822
$ws.Range("T9").Value = 0.01246014137947822

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.65779933333333
$ws.Range("H10").Value = 43.973398
$ws.Range("I10").Value = 0.243637901982001
$ws.Range("J10").Value = 0.243637901982001
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.489401
$ws.Range("N10").Value = 4.468203
$ws.Range("O10").Value = 0.04553544112151264
$ws.Range("P10").Value = 0.04553544112151264
$ws.Range("Q10").Value = 21.831340984866
$ws.Range("R10").Value = 196.482068863794
$ws.Range("S10").Value = 0.01109415934067027
$ws.Range("T10").Value = 0.01109415934067027
